# Add a new "index" column (numeric) as column C on both sheets, shifting the
# existing C:J columns right to D:K, and refresh the "spouse" column
# (Individuals) / "are divorced" column (Families) whose meaning changed
# (they now hold an aggregated list / boolean instead of the old raw value).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: Individuals
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Individuals")

# Insert a new column before the current column C ("name"); this shifts
# name/gender/birthday/child/alive/age/spouse/death from C..J to D..K and
# copies the header style (bold/border) from the neighbouring header cell.
$ws1.Columns("C").Insert()

$ws1.Range("C1").Value = "index"

$ws1.Range("C2").Value = 15
$ws1.Range("C3").Value = 24
$ws1.Range("C4").Value = 35
$ws1.Range("C5").Value = 44
$ws1.Range("C6").Value = 53
$ws1.Range("C7").Value = 62
$ws1.Range("C8").Value = 71
$ws1.Range("C9").Value = 80
$ws1.Range("C10").Value = 88
$ws1.Range("C11").Value = 100
$ws1.Range("C12").Value = 111

# "spouse" (now column J after the shift) becomes the list of family ids the
# individual belongs to as a spouse; refresh the rows whose value changed.
$ws1.Range("J3").Value = "['F1', 'F2']"
$ws1.Range("J4").Value = "['F1']"
$ws1.Range("J6").Value = "['F4']"
$ws1.Range("J7").Value = "['F4']"
$ws1.Range("J8").Value = "['F2']"
$ws1.Range("J10").Value = "['F3', 'F5']"
$ws1.Range("J11").Value = "['F3']"
$ws1.Range("J12").Value = "['F5']"

# ---------------------------------------------------------------------------
# Sheet 2: Families
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Families")

# Insert a new column before the current column C ("Husband ID"); this shifts
# Husband ID/Name, Wife ID/Name, children, married, divorced, are divorced
# from C..J to D..K.
$ws2.Columns("C").Insert()

$ws2.Range("C1").Value = "index"

$ws2.Range("C2").Value = 122
$ws2.Range("C3").Value = 133
$ws2.Range("C4").Value = 141
$ws2.Range("C5").Value = 158
$ws2.Range("C6").Value = 151

# The "divorced"/"are divorced" headers swap order relative to a plain
# column shift: the new "are divorced" boolean column lands in J (with
# "divorced" pushed out to K).
$ws2.Range("J1").Value = "are divorced"
$ws2.Range("K1").Value = "divorced"

# "are divorced" becomes a real boolean column instead of the mostly-empty
# text column it used to be.
$ws2.Range("J2").Value = $true
$ws2.Range("J3").Value = $false
$ws2.Range("J4").Value = $false
$ws2.Range("J5").Value = $false
$ws2.Range("J6").Value = $false
